# Harmonize parameters w Leander thesis
$wb = $excel.ActiveWorkbook

$wsAll    = $wb.Worksheets.Item("All")
$wsLarge  = $wb.Worksheets.Item("Large")
$wsMedium = $wb.Worksheets.Item("Medium")
$wsSmall  = $wb.Worksheets.Item("Small")

# "Medium" sheet: Maximum capacity (GW) 4.7 -> 4.2
$wsMedium.Range("B4").Value = 4.2

# "Small" sheet: Pipeline capex (euros) formula 1.5 * 1000000 -> 0.09 * 1000000
$wsSmall.Range("B2").Formula = "=0.09 * 1000000"

# Update each sheet's selection to match the saved view state, then leave
# "Small" as the active sheet/tab (matches workbookView activeTab + tabSelected).
$wsAll.Activate()
$wsAll.Range("B6").Select()

$wsLarge.Activate()
$wsLarge.Range("B3").Select()

$wsMedium.Activate()
$wsMedium.Range("B5").Select()

$wsSmall.Activate()
$wsSmall.Range("B2").Select()
